$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.445.74"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "3.553.03"
$ws.Range("E3").Value = "  +3.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.60"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.10"
$ws.Range("E6").Value = "  +2.42%  "

$ws.Range("D7").Value = "3.552.92"
$ws.Range("E7").Value = "  +3.06%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  -0.23%  "

$ws.Range("E10").Value = "  +2.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.03"
$ws.Range("E11").Value = "  -5.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.396"
$ws.Range("E12").Value = "  +4.39%  "

$ws.Range("D13").Value = "4.157.33"
$ws.Range("E13").Value = "  +3.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000188"
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.16"
$ws.Range("E15").Value = "  +1.78%  "

$ws.Range("D16").Value = "3.555.58"
$ws.Range("E16").Value = "  +3.32%  "

$ws.Range("E17").Value = "  +1.47%  "

$ws.Range("D18").Value = "65.488.56"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.17"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.90"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.31"
$ws.Range("E21").Value = "  +4.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "395.99"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("E23").Value = "  +4.09%  "

$ws.Range("D24").Value = "3.696.87"
$ws.Range("E24").Value = "  +3.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.23"
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000117"
$ws.Range("E27").Value = "  +9.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.93"
$ws.Range("E28").Value = "  +10.07%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.37"
$ws.Range("E31").Value = "  +1.63%  "

$ws.Range("D32").Value = "3.564.81"
$ws.Range("E32").Value = "  +3.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.29"
$ws.Range("E36").Value = "  +5.78%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.11"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  +2.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.67"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.06"
$ws.Range("E40").Value = "  +5.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0812"
$ws.Range("E41").Value = "  +4.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.833"
$ws.Range("E42").Value = "  +1.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.64"
$ws.Range("E43").Value = "  +15.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.99"
$ws.Range("E44").Value = "  -1.15%  "

$ws.Range("E45").Value = "  +0.11%  "

$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("E47").Value = "  +4.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +8.14%  "

$ws.Range("D49").Value = "2.445.58"
$ws.Range("E49").Value = "  +10.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  +3.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").Value = "  +19.56%  "

# Row 33-35 rotation: Kaspa moves up to row 33, USDe and EthereumClassic shift down
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.148"
$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.84"
$ws.Range("E35").Value = "  +3.91%  "
